$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9720592498779297
$ws.Range("B1").Value = 1.91448438167572
$ws.Range("C1").Value = 5.088549613952637
$ws.Range("D1").Value = 1.64069652557373
$ws.Range("E1").Value = 0.7601715922355652
